$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.007.85'
$ws.Range("E2").Value = '  -3.65%  '
$ws.Range("D3").Value = '1.646.70'
$ws.Range("E3").Value = '  -5.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.97'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -5.57%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4831'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -6.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2600'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -5.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06005'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07195'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("D11").Value = '1.648.09'
$ws.Range("E11").Value = '  -5.35%  '
$ws.Range("E12").Value = '  -2.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6204'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -4.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.513'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '72.85'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -6.06%  '
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9993'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").Value = '25.005.32'
$ws.Range("E18").Value = '  -3.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.41'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.59%  '
$ws.Range("E20").Value = '  -2.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.509'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.64%  '
$ws.Range("D22").Value = '1.855.83'
$ws.Range("E22").Value = '  -5.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.625'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.293'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '131.95'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.86'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.399'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -7.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '102.84'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.671'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.758'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07850'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.580'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04494'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9993'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.594'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9317'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -6.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5822'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.58%  '
$ws.Range("E38").Value = '  -5.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01569'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8477'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +11.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9992'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.820'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '98.19'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3715'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.780'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1149'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.104'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05190'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.80'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9996'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '50.40'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -9.33%  '
